$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits at the end of the "Adding both..."
#    paragraph. In the edited document it has moved to sit inside the new
#    closing line ("I will now drop another object from the sky.") instead.
#    Remove it from its old spot first so we can re-create it later without
#    colliding with the existing id.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2. Locate the blank paragraph that immediately follows the
#    "...daring rescue!" paragraph -- that is where all of the new dialog
#    content gets inserted.
# ---------------------------------------------------------------------------
$rescuePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.TrimEnd() -eq "Now you must use your forces to push the block towards the watchtower for this daring rescue!") {
        $rescuePara = $candidate
        break
    }
}

$targetPara = $rescuePara.Next()

# ---------------------------------------------------------------------------
# 3. Replace that single blank paragraph with the full run of new
#    paragraphs: two blank lines followed by the new lesson dialog, with the
#    relocated "_GoBack" bookmark embedded mid-way through the final line.
# ---------------------------------------------------------------------------
$newBodyXml = @'
<w:p/><w:p/><w:p><w:r><w:t>Excellent! Now that the first block is in place, notice how it took a couple of knights to be able to move the block?</w:t></w:r></w:p><w:p><w:r><w:t>This is because of inertia</w:t></w:r><w:r><w:t>, which causes the block to have more frictional force against the push force of the knights.</w:t></w:r></w:p><w:p><w:r><w:t>Inertia is the tendency of a physical object to resist a change in motion. A change in motion requires acceleration</w:t></w:r><w:r><w:t xml:space="preserve"> due to net force</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>&lt;show highlight on mass&gt;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The mass of an object determines the object’s amount of inertia. The unit of measurement shown </w:t></w:r><w:r><w:t xml:space="preserve">here </w:t></w:r><w:r><w:t xml:space="preserve">is in </w:t></w:r><w:r><w:t>kg (</w:t></w:r><w:r><w:t>kilograms</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t>, which is 1000 grams per 1 kilogram.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">So, the more mass an object has, the more resistance (inertia) it will have from </w:t></w:r><w:r><w:t>change of motion</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>I will now drop another object from the sky.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> This time with less mass.</w:t></w:r><w:r><w:t xml:space="preserve"> Just one more block to rescue our damsel in distress!</w:t></w:r></w:p>
'@

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $newBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetPara.Range.InsertXML($packageXml)

Write-Output "done"
